# Apply annotation updates to the two results tables ("Proteins" and
# "Vaccine" worksheets): table captions are re-worded from "Genes" to
# "Sites", footnote wording is clarified ("positive selection" ->
# "diversifying selection", extra FUBAR/MEME detail, and an added
# sentence about codon 161 on the Vaccine sheet).

$wb = $excel.ActiveWorkbook
$wsProteins = $wb.Worksheets.Item("Proteins")
$wsVaccine  = $wb.Worksheets.Item("Vaccine")

# ---------------------------------------------------------------------
# Table captions first (both sheets), then footnotes (both sheets) -
# this keeps the shared-string table ordered the same way the author's
# edit produced it.
# ---------------------------------------------------------------------

# Table 4 caption (Proteins sheet)
$wsProteins.Range("A1").Value2 = 'Table 4: Sites under Diversifying Selection'

# Table 5 caption (Vaccine sheet) - also fixes the accidental double
# space before "Introduction"
$wsVaccine.Range("A1").Value2 = 'Table 5: Sites under Diversifying Selection after the Introduction of Modern Vaccines'

# Footnote under Table 4 (merged A12:G12, Proteins sheet) - rebuild as
# rich text so the superscripted "+" in "β+" keeps its formatting.
$capNote1 = $wsProteins.Range("A12")
$capNote1.Value2 = 'The number of sequences accepted by HyPhy for analysis is shown by "n". "β-α" is the difference between synonymous (α) and non-synonymous (β) substitution rates over sites and a positive value denotes diversifying selection. "P.Pr." is the posterior probability of diversifying selection. "β+" is the unconstrained estimate for the non-synonymous rate in MEME. All listed proteins were found under both pervasive and episodic diversifying selection by FUBAR (P.Pr.>0.9) and MEME (p-value<0.1) respectively. '
$capNote1.Characters(292, 1).Font.Superscript = $true

# Footnote under Table 5 (merged A6:F6, Vaccine sheet)
$capNote2 = $wsVaccine.Range("A6")
$capNote2.Value2 = 'The number of sequences accepted by HyPhy for analysis is shown by "n". "β-α" is the difference between synonymous (α) and non-synonymous (β) substitution rates over sites and a positive value denotes diversifying selection."P.Pr." is the posterior probability of diversifying selection. "β+" is the unconstrained estimate for the non-synonymous rate in MEME. 
*Only codon 2 was found under both pervasive and episodic diversifying selection by FUBAR (P.Pr.>0.9) and MEME (p-value<0.1) respectively. Codon 161 was found under pervasive diversifying selection by FUBAR (P.Pr.>0.9) but not detected for episodic positive selection by MEME (p-value>0.1).'
$capNote2.Characters(291, 1).Font.Superscript = $true

# Re-select the whole footnote block on each sheet (matches the author
# re-selecting the block after editing it).
$wsProteins.Range("A1:G12").Select()
$wsVaccine.Range("A1:F6").Select()
